$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking price cells as Text so Excel keeps the
# exact original string (trailing zeros, leading zeros, etc.) instead
# of coercing it to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.881.80"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.889.72"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "467.29"
$ws.Range("E5").Value = "  +10.28%  "
$ws.Range("D6").Value = "142.80"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +10.60%  "
$ws.Range("D11").Value = "0.0000332"
$ws.Range("E11").Value = "  +12.07%  "
$ws.Range("D12").Value = "42.84"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "4.507.68"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").Value = "10.30"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "14.90"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "3.893.92"
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "19.78"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "1.16"
$ws.Range("E19").Value = "  +4.90%  "
$ws.Range("D20").Value = "67.103.55"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "429.05"
$ws.Range("E21").Value = "  +7.57%  "
$ws.Range("D22").Value = "14.72"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "3.33"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "88.08"
$ws.Range("E24").Value = "  +5.32%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "3.54"
$ws.Range("E25").Value = "  +9.71%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "38.35"
$ws.Range("E26").Value = "  +6.00%  "
$ws.Range("D27").Value = "5.76"
$ws.Range("E27").Value = "  +6.42%  "
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "730.37"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("D31").Value = "13.73"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "43.15"
$ws.Range("E34").Value = "  +7.27%  "
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("D36").Value = "57.26"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0784"
$ws.Range("E38").Value = "  +21.90%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "5.39"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").Value = "0.0474"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("E41").Value = "  +8.39%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.338"
$ws.Range("E42").Value = "  +7.42%  "
$ws.Range("D43").Value = "0.141"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D47").Value = "2.17"
$ws.Range("E47").Value = "  +6.98%  "
$ws.Range("D48").Value = "3.39"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").Value = "144.52"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "2.87"
$ws.Range("E51").Value = "  +4.94%  "
